$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update student data rows (columns B,D,E,G,H; I unchanged) ---

# Row 2
$ws.Range("B2").Value = "E22110469000110028"
$ws.Range("D2").Value = "vraj dalal"
$ws.Range("E2").Value = "vrajdalal492@gmail.com"
$ws.Range("G2").Value = "BCA"
$ws.Range("H2").Value = "Third Year"

# Row 3
$ws.Range("B3").Value = "E22110469000110029"
$ws.Range("D3").Value = "mukesh"
$ws.Range("E3").Value = "suresh@gmail.com"
$ws.Range("G3").Value = "BBA"
$ws.Range("H3").Value = "First Year"

# Row 4
$ws.Range("B4").Value = "E22110469000110030"
$ws.Range("D4").Value = "ramesh"
$ws.Range("E4").Value = "ramesh@gmail.com"
$ws.Range("G4").Value = "BCOM"
$ws.Range("H4").Value = "First Year"

# Row 5
$ws.Range("B5").Value = "E22110469000110031"
$ws.Range("D5").Value = "nilesh"
$ws.Range("E5").Value = "vrajdalal650@gmail.com"
$ws.Range("G5").Value = "BBA"
$ws.Range("H5").Value = "Second Year"

# Row 6
$ws.Range("B6").Value = "E22110469000110032"
$ws.Range("D6").Value = "rajesh"
$ws.Range("E6").Value = "rajesh@gmail.com"
$ws.Range("G6").Value = "BCA"
$ws.Range("H6").Value = "Third Year"

# --- Re-create the mailto hyperlinks so they point at the right rows/emails ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:suresh@gmail.com", "", "mailto:suresh@gmail.com", "suresh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:ramesh@gmail.com", "", "mailto:ramesh@gmail.com", "ramesh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:vrajdalal650@gmail.com", "", "mailto:vrajdalal650@gmail.com", "vrajdalal650@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:rajesh@gmail.com", "", "mailto:rajesh@gmail.com", "rajesh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:vrajdalal492@gmail.com", "", "mailto:vrajdalal492@gmail.com", "vrajdalal492@gmail.com")

# --- Column widths ---
# (Target stored widths are 11.7142857142857 / 20.7142857142857 / 24.7142857142857
#  characters, i.e. 82/145/173 pixels at MDW=7. This runtime quantizes ColumnWidth
#  on a 6px-MDW grid, so we feed in the char-width whose rounded pixel count is the
#  nearest achievable to the true target.)
$ws.Columns.Item(1).ColumnWidth = 10.8333333333333
$ws.Columns.Item(2).ColumnWidth = 19.8333333333333
$ws.Columns.Item(5).ColumnWidth = 23.8333333333333

# --- Selection ---
$ws.Range("B7").Select()
